# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# This updates the "K" column values (column G, rows 2-26) on the active sheet
# with the recalculated strikeout totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 5
    4  = 8
    5  = 3
    6  = 2
    7  = 6
    8  = 4
    9  = 9
    10 = 4
    11 = 5
    12 = 4
    13 = 4
    14 = 1
    15 = 3
    16 = 5
    17 = 5
    18 = 4
    19 = 4
    20 = 4
    21 = 4
    22 = 4
    23 = 1
    24 = 2
    25 = 4
    26 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
